$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$cs = $nm.ColorScheme
$cs.Colors(1).RGB = 999999
